$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new student rows (14-16), following the same pattern as
# the preceding rows: column B = "power poite", C = "dezembro de 2006",
# D = "dezembro de 2007", E incrementing from the last value (110).

$ws.Range("A14").Value = "ana maria dos santos"
$ws.Range("B14").Value = "power poite"
$ws.Range("C14").Value = "dezembro de 2006"
$ws.Range("D14").Value = "dezembro de 2007"
$ws.Range("E14").Value = 111

$ws.Range("A15").Value = "isa gabrielly de oliveira"
$ws.Range("B15").Value = "power poite"
$ws.Range("C15").Value = "dezembro de 2006"
$ws.Range("D15").Value = "dezembro de 2007"
$ws.Range("E15").Value = 112

$ws.Range("A16").Value = "seu maico de aumenda"
$ws.Range("B16").Value = "power poite"
$ws.Range("C16").Value = "dezembro de 2006"
$ws.Range("D16").Value = "dezembro de 2007"
$ws.Range("E16").Value = 113

# Leave the cursor resting a couple rows below the new data, matching
# where the author ended up after typing the last entry.
$null = $ws.Range("A19").Select()
